$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.096773137147352
$ws.Range("C2").Value = 0.1539609373135704
$ws.Range("D2").Value = 0.1578897478055623
$ws.Range("F2").Value = 1.51453054433793
$ws.Range("G2").Value = 0.00246849754527558
$ws.Range("I2").Value = 0.9084574317527583
$ws.Range("J2").Value = 0.183131137133727
$ws.Range("L2").Value = 0.3672861355196133
$ws.Range("M2").Value = 0.3026410078799415
$ws.Range("N2").Value = 1.549546798466331
$ws.Range("O2").Value = 3.751721784098265
$ws.Range("B3").Value = 1.019581426637899
$ws.Range("C3").Value = 0.1417995896655668
$ws.Range("D3").Value = 0.157844914693321
$ws.Range("F3").Value = 1.518461065220563
$ws.Range("G3").Value = 0.002471202567580878
$ws.Range("I3").Value = 0.9159447299479808
$ws.Range("J3").Value = 0.1842172475221937
$ws.Range("L3").Value = 0.364461574930921
$ws.Range("M3").Value = 0.2899106986700275
$ws.Range("N3").Value = 1.56225450252002
$ws.Range("O3").Value = 3.762499558357945
$ws.Range("B4").Value = 0.972402053314255
$ws.Range("C4").Value = 0.1342750391756766
$ws.Range("D4").Value = 0.1578573912103813
$ws.Range("F4").Value = 1.521638764243704
$ws.Range("G4").Value = 0.00247295405883674
$ws.Range("I4").Value = 0.9209900786908598
$ws.Range("J4").Value = 0.1849206735593469
$ws.Range("L4").Value = 0.362856146443491
$ws.Range("M4").Value = 0.2821854883637442
$ws.Range("N4").Value = 1.570584180868352
$ws.Range("O4").Value = 3.771149280661717
$ws.Range("B5").Value = 0.953231986318599
$ws.Range("C5").Value = 0.1311944232397622
$ws.Range("D5").Value = 0.1578725727828179
$ws.Range("F5").Value = 1.523126026618989
$ws.Range("G5").Value = 0.002473690654212915
$ws.Range("I5").Value = 0.9231587920337887
$ws.Range("J5").Value = 0.1852165328098456
$ws.Range("L5").Value = 0.3622344520923448
$ws.Range("M5").Value = 0.2790606151399828
$ws.Range("N5").Value = 1.574111255064992
$ws.Range("O5").Value = 3.775185225898724
$ws.Range("B6").Value = 0.950052230360285
$ws.Range("C6").Value = 0.1306820297694173
$ws.Range("D6").Value = 0.1578757046378882
$ws.Range("F6").Value = 1.523384604968214
$ws.Range("G6").Value = 0.002473814347428035
$ws.Range("I6").Value = 0.9235257119003144
$ws.Range("J6").Value = 0.1852662166097043
$ws.Range("L6").Value = 0.3621331887405148
$ws.Range("M6").Value = 0.2785431421710172
$ws.Range("N6").Value = 1.574704939136581
$ws.Range("O6").Value = 3.775886265045386
$ws.Range("B7").Value = 0.9721432908072529
$ws.Range("C7").Value = 0.1342335506258081
$ws.Range("D7").Value = 0.1578575550195573
$ws.Range("F7").Value = 1.521658043107195
$ws.Range("G7").Value = 0.002472963900289555
$ws.Range("I7").Value = 0.9210188704159599
$ws.Range("J7").Value = 0.1849246263177271
$ws.Range("L7").Value = 0.3628476301702719
$ws.Range("M7").Value = 0.2821432508893409
$ws.Range("N7").Value = 1.570631210917696
$ws.Range("O7").Value = 3.771201641198843
$ws.Range("B8").Value = 1.070113303811297
$ws.Range("C8").Value = 0.1497797066108149
$ws.Range("D8").Value = 0.1578660089452981
$ws.Range("F8").Value = 1.515727220004457
$ws.Range("G8").Value = 0.002469411474126048
$ws.Range("I8").Value = 0.9109460512667624
$ws.Range("J8").Value = 0.1834980521680434
$ws.Range("L8").Value = 0.3662855561311176
$ws.Range("M8").Value = 0.2982328137322057
$ws.Range("N8").Value = 1.553819087656365
$ws.Range("O8").Value = 3.755016257005337
$ws.Range("B9").Value = 1.263898127792743
$ws.Range("C9").Value = 0.1798050731755723
$ws.Range("D9").Value = 0.1581984354951231
$ws.Range("F9").Value = 1.510157527044797
$ws.Range("G9").Value = 0.002463160955284853
$ws.Range("I9").Value = 0.894748634576402
$ws.Range("J9").Value = 0.1809897744726205
$ws.Range("L9").Value = 0.3740448462022954
$ws.Range("M9").Value = 0.3304984600846339
$ws.Range("N9").Value = 1.525027463618777
$ws.Range("O9").Value = 3.739399178452118
$ws.Range("B10").Value = 1.407230593092322
$ws.Range("C10").Value = 0.2015793245740554
$ws.Range("D10").Value = 0.1586332563138981
$ws.Range("F10").Value = 1.509756709839039
$ws.Range("G10").Value = 0.002459000764525984
$ws.Range("I10").Value = 0.8850157517982495
$ws.Range("J10").Value = 0.179322154455904
$ws.Range("L10").Value = 0.38036020872714
$ws.Range("M10").Value = 0.354628518854085
$ws.Range("N10").Value = 1.50641321187188
$ws.Range("O10").Value = 3.737757202606076
$ws.Range("B11").Value = 1.472633263151067
$ws.Range("C11").Value = 0.2114221689823239
$ws.Range("D11").Value = 0.1588720463517177
$ws.Range("F11").Value = 1.510375289588595
$ws.Range("G11").Value = 0.002457201092715891
$ws.Range("I11").Value = 0.881058652246054
$ws.Range("J11").Value = 0.1786013196243648
$ws.Range("L11").Value = 0.3833655406031795
$ws.Range("M11").Value = 0.3656959449676052
$ws.Range("N11").Value = 1.4984949210289
$ws.Range("O11").Value = 3.739146207142568
$ws.Range("B12").Value = 1.497427051377088
$ws.Range("C12").Value = 0.2151403101790095
$ws.Range("D12").Value = 0.1589683278128149
$ws.Range("F12").Value = 1.510724616142127
$ws.Range("G12").Value = 0.002456532882311453
$ws.Range("I12").Value = 0.8796278537702733
$ws.Range("J12").Value = 0.1783337729366918
$ws.Range("L12").Value = 0.3845225115056081
$ws.Range("M12").Value = 0.3698996747289129
$ws.Range("N12").Value = 1.495575378504533
$ws.Range("O12").Value = 3.739979331112664
$ws.Range("B13").Value = 1.49208607958326
$ws.Range("C13").Value = 0.21433995059013
$ws.Range("D13").Value = 0.1589473319828301
$ws.Range("F13").Value = 1.510644265417454
$ws.Range("G13").Value = 0.002456676203367732
$ws.Range("I13").Value = 0.8799329918813257
$ws.Range("J13").Value = 0.1783911531658982
$ws.Range("L13").Value = 0.384272497821101
$ws.Range("M13").Value = 0.3689937644065751
$ws.Range("N13").Value = 1.496200644224068
$ws.Range("O13").Value = 3.739786243061133
$ws.Range("B14").Value = 1.474672526414793
$ws.Range("C14").Value = 0.2117282466072368
$ws.Range("D14").Value = 0.1588798503185629
$ws.Range("F14").Value = 1.510401722940429
$ws.Range("G14").Value = 0.002457145852573758
$ws.Range("I14").Value = 0.8809395832820428
$ws.Range("J14").Value = 0.1785791998957276
$ws.Range("L14").Value = 0.3834603469042719
$ws.Range("M14").Value = 0.366041534728744
$ws.Range("N14").Value = 1.498253146740232
$ws.Range("O14").Value = 3.739208593302351
$ws.Range("B15").Value = 1.464009720862123
$ws.Range("C15").Value = 0.2101273086735773
$ws.Range("D15").Value = 0.1588392774857041
$ws.Range("F15").Value = 1.5102681434049
$ws.Range("G15").Value = 0.002457435254764234
$ws.Range("I15").Value = 0.881564962733421
$ws.Range("J15").Value = 0.1786950890939476
$ws.Range("L15").Value = 0.3829653407659208
$ws.Range("M15").Value = 0.3642348589450393
$ws.Range("N15").Value = 1.499520642580357
$ws.Range("O15").Value = 3.738894764312732
$ws.Range("B16").Value = 1.402960267431354
$ws.Range("C16").Value = 0.2009348026668647
$ws.Range("D16").Value = 0.158618472425502
$ws.Range("F16").Value = 1.50973239323541
$ws.Range("G16").Value = 0.002459120239874947
$ws.Range("I16").Value = 0.8852838257522819
$ws.Range("J16").Value = 0.1793700216131899
$ws.Range("L16").Value = 0.3801664573405645
$ws.Range("M16").Value = 0.3539070341637043
$ws.Range("N16").Value = 1.506941741055037
$ws.Range("O16").Value = 3.737709407919453
$ws.Range("B17").Value = 1.365558533179978
$ws.Range("C17").Value = 0.195279408228572
$ws.Range("D17").Value = 0.1584934871761732
$ws.Range("F17").Value = 1.509608799570955
$ws.Range("G17").Value = 0.002460177653065234
$ws.Range("I17").Value = 0.8876857329692953
$ws.Range("J17").Value = 0.1797937356309349
$ws.Range("L17").Value = 0.3784832633842399
$ws.Range("M17").Value = 0.3475942357759294
$ws.Range("N17").Value = 1.511635018619003
$ws.Range("O17").Value = 3.737529291951518
$ws.Range("B18").Value = 1.344064963239362
$ws.Range("C18").Value = 0.192020713846432
$ws.Range("D18").Value = 0.1584254591531362
$ws.Range("F18").Value = 1.509613102731393
$ws.Range("G18").Value = 0.00246079458957754
$ws.Range("I18").Value = 0.8891115250342523
$ws.Range("J18").Value = 0.1800410010130253
$ws.Range("L18").Value = 0.3775276084303272
$ws.Range("M18").Value = 0.3439718202485409
$ws.Range("N18").Value = 1.514386185870535
$ws.Range("O18").Value = 3.73762674009285
$ws.Range("B19").Value = 1.336790906130886
$ws.Range("C19").Value = 0.1909163741909765
$ws.Range("D19").Value = 0.158403090157563
$ws.Range("F19").Value = 1.509627510811853
$ws.Range("G19").Value = 0.002461004976838375
$ws.Range("I19").Value = 0.8896018779809864
$ws.Range("J19").Value = 0.180125332047413
$ws.Range("L19").Value = 0.3772061863556218
$ws.Range("M19").Value = 0.342746807858866
$ws.Range("N19").Value = 1.515326568497827
$ws.Range("O19").Value = 3.737694264463357
$ws.Range("B20").Value = 1.369538064858318
$ws.Range("C20").Value = 0.1958820423906786
$ws.Range("D20").Value = 0.1585063928239592
$ws.Range("F20").Value = 1.509614154269869
$ws.Range("G20").Value = 0.002460064185515257
$ws.Range("I20").Value = 0.8874254630285634
$ws.Range("J20").Value = 0.1797482625903672
$ws.Range("L20").Value = 0.3786611522025254
$ws.Range("M20").Value = 0.3482653622309684
$ws.Range("N20").Value = 1.511130059055255
$ws.Range("O20").Value = 3.737527658141943
$ws.Range("B21").Value = 1.479786582746101
$ws.Range("C21").Value = 0.2124956164566072
$ws.Range("D21").Value = 0.1588995126609518
$ws.Range("F21").Value = 1.510469840895354
$ws.Range("G21").Value = 0.002457007545152814
$ws.Range("I21").Value = 0.8806420862969944
$ws.Range("J21").Value = 0.1785238190860126
$ws.Range("L21").Value = 0.3836983829344547
$ws.Range("M21").Value = 0.3669083327674656
$ws.Range("N21").Value = 1.497648135017883
$ws.Range("O21").Value = 3.739369927474769
$ws.Range("B22").Value = 1.55199828543681
$ws.Range("C22").Value = 0.2233002539122992
$ws.Range("D22").Value = 0.1591905602816013
$ws.Range("F22").Value = 1.511699837709671
$ws.Range("G22").Value = 0.002455087268075196
$ws.Range("I22").Value = 0.8766031953391078
$ws.Range("J22").Value = 0.1777551472598073
$ws.Range("L22").Value = 0.3871007075607906
$ws.Range("M22").Value = 0.3791666934745521
$ws.Range("N22").Value = 1.489297027768778
$ws.Range("O22").Value = 3.742364187856253
$ws.Range("B23").Value = 1.513443600706637
$ws.Range("C23").Value = 0.2175385452325997
$ws.Range("D23").Value = 0.1590321128504044
$ws.Range("F23").Value = 1.510982019744318
$ws.Range("G23").Value = 0.002456105092827154
$ws.Range("I23").Value = 0.8787227294516384
$ws.Range("J23").Value = 0.1781625175640187
$ws.Range("L23").Value = 0.385274781532317
$ws.Range("M23").Value = 0.3726174891138214
$ws.Range("N23").Value = 1.493712091069632
$ws.Range("O23").Value = 3.740602291706352
$ws.Range("B24").Value = 1.367738889856525
$ws.Range("C24").Value = 0.1956096143170498
$ws.Range("D24").Value = 0.1585005462558371
$ws.Range("F24").Value = 1.509611498665777
$ws.Range("G24").Value = 0.002460115456156165
$ws.Range("I24").Value = 0.887542991197563
$ws.Range("J24").Value = 0.1797688095297643
$ws.Range("L24").Value = 0.3785806911690344
$ws.Range("M24").Value = 0.3479619244576
$ws.Range("N24").Value = 1.511358186364284
$ws.Range("O24").Value = 3.737527770683755
$ws.Range("B25").Value = 1.211301734018207
$ws.Range("C25").Value = 0.1717322012089539
$ws.Range("D25").Value = 0.1580748949221729
$ws.Range("F25").Value = 1.511015804109775
$ws.Range("G25").Value = 0.002464775702763682
$ws.Range("I25").Value = 0.8987498490143935
$ws.Range("J25").Value = 0.1816374813314994
$ws.Range("L25").Value = 0.3718374093605235
$ws.Range("M25").Value = 0.3216944228817056
$ws.Range("N25").Value = 1.532369975524915
$ws.Range("O25").Value = 3.741897659964224
